$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jobs")
$summary = $wb.Worksheets.Item("Summary")

# Widen column F (JobDescriptionSummary) to match column G (width 80)
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Fill in JobDescriptionSummary (column F) text that failed to scrape previously
$f2 = @"
Client Engagement & Consulting
Interact directly with clients and stakeholders to understand complex business problems.
Translate business needs into AI/ML opportunities and articulate solution strategies clearly.
Present solution architectures, effort estimates, timelines, and potential ROI.
Solution Design & Architecture
Define end-to-end AI architecture including data pipelines, model selection, AI engines, integration layers, APIs, and deployment strategies.
Make informed choices of AI/ML/DL algorithms, frameworks, and platforms based on problem type and constraints.
Ensure the architecture is scalable, modular, secure, and optimized for performance.
Hands-On Development
Develop critical components or proof-of-concepts when needed.
Perform hands-on design and coding for pipelines, model training, inferencing, or integrations.
Create reusable components and automation scripts for repeatability.
Team Leadership & Delivery Oversight
Build and mentor a high-performing AI development team.
Drive agile development cycles, ensuring timely and quality delivery.
Collaborate with DevOps, Product Managers, and QA teams throughout SDLC.
Client Engagement & Consulting
Interact directly with clients and stakeholders to understand complex business problems.
Translate business needs into AI/ML opportunities and articulate solution strategies clearly.
Present solution architectures, effort estimates, timelines, and potential ROI.
Solution Design & Architecture
Define end-to-end AI architecture including data pipelines, model selection, AI engines, integration layers, APIs, and deployment strategies.
Make informed choices of AI/ML/DL algorithms, frameworks, and platforms based on problem type and constraints.
Ensure the architecture is scalable, modular, secure, and optimized for performance.
Hands-On Development
Develop critical components or proof-of-concepts when needed.
Perform hands-on design and coding for pipelines, model training, inferencing, or integrations.
Create reusable components and automation scripts for repeatability.
Team Leadership & Delivery Oversight
Build and mentor a high-performing AI development team.
Drive agile development cycles, ensuring timely and quality delivery.
Collaborate with DevOps, Product Managers, and QA teams throughout SDLC.
"@
$ws.Range("F2").Value = $f2
$ws.Rows.Item(2).RowHeight = 60

$f4 = @"
Develop ML and Deep Learning SolutionsCreate predictive, classification, and optimization models using supervised, unsupervised, and reinforcement learning.
Design and Implement Generative AI SystemsBuild and fine-tune Large Language Models (LLMs) and diffusion models for a range of use cases involving text, code, and multi-modal data.
Build Scalable Recommendation EnginesDevelop and optimize recommendation systems using collaborative filtering, content-based filtering, hybrid models, or sequential models.
Cloud-Native ML EngineeringDeploy and manage machine learning pipelines and APIs in cloud environments (AWS, GCP, Azure),ensuring scalability, observability, and cost efficiency.
End-to-End ML Lifecycle OwnershipOwn the model lifecycle from feature engineering and experimentation to deployment, CI/CD, monitoring, and iteration.
Collaborate Across FunctionsWork with cross-functional teams including product managers, engineers, and data scientists to translate business goals into AI-powered solutions.
"@
$ws.Range("F4").Value = $f4
$ws.Rows.Item(4).RowHeight = 60

$f5 = @"
Data mining or extracting usable data from valuable data sources
Using machine learning tools to select features, create and optimize classifiers
Carrying out the preprocessing of structured and unstructured data
Enhancing data collection procedures to include all relevant information for developing analytic systems
Processing, cleansing, and validating the integrity of data to be used for analysis
Analyzing large amounts of information to find patterns and solutions
Developing prediction systems and machine learning algorithms
Presenting results in a clear manner
Propose solutions and strategies to tackle business challenges
Collaborate with Business and IT teams
"@
$ws.Range("F5").Value = $f5
$ws.Rows.Item(5).RowHeight = 60

$f9 = @"
The position holder is primarily responsible for revenue generation. through acquisition of new customers and strategic alliances with existing customers.
Research and identify new business opportunities – including new markets, growth areas, trends, customers, partnerships, products and services.
Implement effective Sales strategies and contribute to revenue targets.
Develop a growth strategy focused both on financial gain and customer satisfaction.
Understand client needs and present appropriate IT solutions to meet their requirements.
Build strategic alliances, joint ventures or collaborations that create mutually beneficial opportunities.
Prepare and deliver sales proposals, quotes, and contracts to prospective clients.
Generate regular reports on sales performance and forecasts.
Maintain accurate records of sales activities, including sales calls, presentations, and closed deals.
Work closely with internal teams to develop customized solutions based on client needs.
Build and maintain strong, long-lasting relationships with clients.
"@
$ws.Range("F9").Value = $f9
$ws.Rows.Item(9).RowHeight = 60

$f10 = @"
Collaborate closely with the business team to comprehend project-specific business requirements, fostering effective relationships with clients.
Define project deliverables, resource needs, and a comprehensive work plan, overseeing their development and successful delivery.
Guarantee adherence to financial parameters outlined by the business unit throughout the project.
Coordinate and manage project teams, providing guidance, support, and motivation.
Engage with stakeholders to gather requirements, provide updates, and manage expectations.
Oversee the implementation of technical solutions, including code reviews, testing, and deployment.
Develop accurate project estimates and comprehensive plans, incorporating contingency measures from project initiation to completion.
Collaborate with the QA team to ensure projects meet stringent quality standards and align with all Techversant processes.
Effectively manage, monitor, and inspire the assigned project team to achieve optimal performance.
Conduct compelling project presentations and compile comprehensive reports for stakeholders.
Facilitate management training sessions to acquaint teams with organizational processes and procedures.
Proactively identify and report potential risks and issues to management for timely resolution.
"@
$ws.Range("F10").Value = $f10
$ws.Rows.Item(10).RowHeight = 60

$f12 = @"
Develop and maintain web applications using server-side technologies such as PHP or Python.
Collaborate with the UI/UX team to implement responsive and user-friendly web interfaces using HTML, CSS, and JavaScript.
Write clean, efficient, and maintainable code, following OOP (Object-Oriented Programming) principles.
Troubleshoot, debug, and optimize web applications to ensure maximum performance and usability.
Assist in designing and implementing new features and functionalities.
Ensure responsive design and cross-browser compatibility for a seamless user experience across devices.
Participate in code reviews, providing constructive feedback and learning from senior developers.
Stay updated with industry trends and continuously improve your knowledge of new technologies and frameworks.
"@
$ws.Range("F12").Value = $f12
$ws.Rows.Item(12).RowHeight = 60

$f13 = @"
Design, implement and support the technical solution.
Participate actively in all phases of the application development lifecycle.
Work with user experience designers to ensure all user interactions are
Implemented correctly and optimized for performance so that we can build an amazing user experience.
Work within an agile group for varying degrees to deliver solutions.
You will be responsible for writing server-side web application logic in JavaScript and/or variants of it.
"@
$ws.Range("F13").Value = $f13
$ws.Rows.Item(13).RowHeight = 60

$f14 = @"
Lead a team of developers and provide technical guidance, mentoring, and performance management.
Lead code reviews to ensure code quality, consistency, and adherence to coding standards.
Design and develop high-quality software solutions using .NET Core, Angular or React, MS SQl and Azure technologies
Collaborate with cross-functional teams to understand project requirements and translate them into technical specifications.
Stay up-to-date with the latest .NET technologies, frameworks, and best practices.
Troubleshoot and resolve complex technical issues across the application stack.
Work closely with product managers and designers to ensure successful project delivery.
Manage and prioritize project tasks, ensuring deadlines are met and milestones achieved.
Architect scalable and maintainable solutions that align with the company’s technical vision.
Participate in project planning, estimation, and documentation processes.
Should take effort to understand the code in other technology project and convert the project into dotnet project, if such requirement comes.
Ready to learn new technologies.
Ensures code reviews, unit testing and integration testing is completed efficiently and effectively.
"@
$ws.Range("F14").Value = $f14
$ws.Rows.Item(14).RowHeight = 60

$f16 = @"
Work closely with team members to support their use of the platform to deliver complete solutions to customers
Participate in all phases of an agile development cycle – plan, design, implement, review, test, deploy, document, and training
Troubleshoot and resolve customer issues
Strive to make continuous improvement
Architect highly scalable, highly available systems with minimal guidance.
Translate business requirements into scalable and extensible design.
Mentor junior engineers on design, coding and troubleshooting.
Practice and advocate Agile processes in the team
"@
$ws.Range("F16").Value = $f16
$ws.Rows.Item(16).RowHeight = 60

$f18 = @"
Lead the development and implementation of mobile applications using React Native.
Collaborate with cross-functional teams to gather and analyze project requirements.
Previous ownership of app development projects is required.
Develop automated tests to guarantee the product’s top-notch quality.
Stay updated with the latest advancements and changes in mobile platforms and technologies.
Provide technical supervision and support to mentor the team effectively.
Conduct code reviews to ensure adherence to coding standards, best practices, and maintain code quality.
Ensure the application is published and maintained in both the App Store and Google Play Store, adhering to their respective policies and guideline.
Participate in agile development processes, including sprint planning, estimation, and retrospectives.
"@
$ws.Range("F18").Value = $f18
$ws.Rows.Item(18).RowHeight = 60

$f19 = @"
The position holder is primarily responsible for revenue generation. through acquisition of new customers and strategic alliances with existing customers.
Research and identify new business opportunities – including new markets, growth areas, trends, customers, partnerships, products and services.
Implement effective Sales strategies and contribute to revenue targets.
Develop a growth strategy focused both on financial gain and customer satisfaction.
Understand client needs and present appropriate IT solutions to meet their requirements.
Build strategic alliances, joint ventures or collaborations that create mutually beneficial opportunities.
Prepare and deliver sales proposals, quotes, and contracts to prospective clients.
Generate regular reports on sales performance and forecasts.
Maintain accurate records of sales activities, including sales calls, presentations, and closed deals.
Work closely with internal teams to develop customized solutions based on client needs.
Build and maintain strong, long-lasting relationships with clients.
Build, lead, mentor, and motivate the business team to achieve sales targets and objectives.
"@
$ws.Range("F19").Value = $f19
$ws.Rows.Item(19).RowHeight = 60

$f20 = @"
Collaborate with cross-functional teams to gather and analyze project requirements, and then translate them into technical specifications and development plans.
Design and develop responsive and intuitive frontend interfaces using Angular, ensuring seamless user experiences across different devices.
Create and maintain robust backend APIs and services using Node.js, Express, and related technologies.
Implement secure authentication and authorization mechanisms for both frontend and backend components.
Write clean, modular, and efficient code that follows coding standards and best practices.
Debug and troubleshoot complex issues across the application stack, ensuring optimal performance and reliability.
Integrate third-party APIs and external services into the applications as needed.
Work closely with UI/UX Designers to implement visually appealing and user-centric interfaces.
Collaborate with other developers to conduct code reviews and maintain code quality.
Mentor and assist junior developers in their technical growth and skill development.
Stay updated with emerging trends and technologies in full stack development and share insights with the team.
Experience in relational and non-relation databases and ability to write complex queries.
Ability to understand business requirement and model it as database entities.
"@
$ws.Range("F20").Value = $f20
$ws.Rows.Item(20).RowHeight = 60

$f21 = @"
Experience in Design, Development, and Delivery of large-scale Web Applications.
Should be competent in code optimization techniques, coding practices and Static code analysis tools.
Should be competent in understanding complex business requirements.
Demonstrates the ability to solve problems and create solutions independently.
Addresses the problems of systems integration, compatibility, and multiple platforms.
Capable enough to communicate with clients to understand the technical requirements and present technical solutions and system design.
Ability to analyse the business requirement, prepare work breakdown structure and provide technical estimate.
Provides good visibility and update on the progress of all activities and tasks, escalating issues and in a timely manner.
Maintain a cohesive spirit to accomplish the team responsibilities.
Proficient enough to be a lead contributor on assigned projects.
Instructs, assigns, directs, and checks the work of other software developers on development team when required.
Should be able to conduct peer reviews on code quality.
Need to ensure the Information Security in all perspectives.
"@
$ws.Range("F21").Value = $f21
$ws.Rows.Item(21).RowHeight = 60

$f22 = @"
Mentor junior members of the team.
Provide timely code reviews for peers.
Develop high-quality software solutions using .NET Core, Angular or React, MS SQl and Azure technologies
Collaborate with cross-functional teams to understand project requirements and translate them into technical specifications.
Stay up-to-date with the latest .NET technologies, frameworks, and best practices.
Troubleshoot and resolve complex technical issues across the application stack.
Actively participates in scrum meetings and analyzes requirements to understand technical/functional requirements. Also possess a capability to discuss requirement with client directly.
Ability to work as an individual contributor.
Ability to assign task to junior members of the team and guide them in their development
Manage and prioritize project tasks after discussing with client, ensuring deadlines are met and milestones achieved.
Participate in project planning, estimation, and documentation processes.
Should take effort to understand the code in other technology project and convert the project into dotnet project, if such requirement comes.
Ready to learn new technologies.
Should have a commitment to complete the task on time without a follow up from senior.
Learn technology business domain and system domain as recommended by the project
Code ,debug,test , document and communicate product/component/feature development stage
Ensures code reviews, unit testing and integration testing is completed efficiently and effectively.
"@
$ws.Range("F22").Value = $f22
$ws.Rows.Item(22).RowHeight = 60

# Update ScrapedDate (column M) for every job row to the new scrape timestamp
$newScrapedDate = "2025-08-03 21:19:22"
for ($row = 2; $row -le 23; $row++) {
    $ws.Range("M$row").Value = $newScrapedDate
}

# Update the Summary sheet Scraping Date (B3) to match
$summary.Range("B3").Value = $newScrapedDate
